# adding buffer description and fraxinus genus
#
# The watchlist table (A:scientific.name, B:status, C:in.anp) gets a new
# row for the bare "Fraxinus" genus, inserted at row 136 (just above
# "Geum rivale"), pushing everything below it down by one row. The
# trailing blank "buffer" row (previously A189, hyperlink-style s="2")
# shifts down to A190 along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 136, shifting rows 136:189 down to 137:190.
[void]$ws.Rows.Item(136).Insert()

# Populate the new row with the Fraxinus genus entry (same status/in.anp
# as its neighbours in the "rare native" / "P" block).
$ws.Range("A136").Value = "Fraxinus"
$ws.Range("B136").Value = "rare native"
$ws.Range("C136").Value = "P"

# Reflect where the author was working when they made the edit.
[void]$ws.Range("D136").Select()

Write-Output "Inserted Fraxinus row at A136"
